# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice / profit calculations to Garuda_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart / Roof Tile
$ws.Range("H19").Value2 = 6466.6665
$ws.Range("I19").Value2 = 7933.3335
$ws.Range("J19").Value2 = 3533.3333
$ws.Range("K19").Value2 = 7933.3335
$ws.Range("L19").Value2 = 3533.3333
$ws.Range("M19").Value2 = -7758.3335
$ws.Range("N19").Value2 = -3883.3333

# Row 53: No Accounting for Waste / Enchanted Electrum Ink
$ws.Range("H53").Value2 = 284.1
$ws.Range("I53").Value2 = 69.888885
$ws.Range("K53").Value2 = 69.888885
$ws.Range("M53").Value2 = 567.111115

# Row 103: Let Loose the Juice / Persimmon Tannin
$ws.Range("H103").Value2 = 538.25
$ws.Range("J103").Value2 = 702.5
$ws.Range("L103").Value2 = 2107.5
$ws.Range("N103").Value2 = -3279.5

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value2 = 1744
$ws.Range("I116").Value2 = 1684.9231
$ws.Range("J116").Value2 = 2000
$ws.Range("K116").Value2 = 1684.9231
$ws.Range("L116").Value2 = 2000
$ws.Range("M116").Value2 = 1757.0769
$ws.Range("N116").Value2 = -8884

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value2 = 1549.1316
$ws.Range("I137").Value2 = 1218.7916
$ws.Range("J137").Value2 = 2115.4285
$ws.Range("K137").Value2 = 3656.3748
$ws.Range("L137").Value2 = 6346.2855
$ws.Range("M137").Value2 = -1106.3748
$ws.Range("N137").Value2 = -11446.2855

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value2 = 1152.49
$ws.Range("I32").Value2 = 1170.6083
$ws.Range("J32").Value2 = 566.6667
$ws.Range("K32").Value2 = 1170.6083
$ws.Range("L32").Value2 = 566.6667
$ws.Range("M32").Value2 = -883.6083000000001
$ws.Range("N32").Value2 = -1140.6667

# Row 42: Kitty Get Your Helm / Steel Elmo
$ws.Range("H42").Value2 = 5000
$ws.Range("I42").Value2 = 5000
$ws.Range("J42").Value2 = 0
$ws.Range("K42").Value2 = 5000
$ws.Range("L42").Value2 = 0
$ws.Range("M42").Value2 = -4514
$ws.Range("N42").ClearContents()

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value2 = 1285
$ws.Range("I45").Value2 = 1121.1
$ws.Range("J45").Value2 = 1394.2667
$ws.Range("K45").Value2 = 1121.1
$ws.Range("L45").Value2 = 1394.2667
$ws.Range("M45").Value2 = -744.0999999999999
$ws.Range("N45").Value2 = -2148.2667

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value2 = 1247.4722
$ws.Range("I61").Value2 = 837.43634
$ws.Range("J61").Value2 = 2574.0588
$ws.Range("K61").Value2 = 837.43634
$ws.Range("L61").Value2 = 2574.0588
$ws.Range("M61").Value2 = -625.43634
$ws.Range("N61").Value2 = -2998.0588

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value2 = 1412.6
$ws.Range("I110").Value2 = 1623.75
$ws.Range("J110").Value2 = 1171.2858
$ws.Range("K110").Value2 = 1623.75
$ws.Range("L110").Value2 = 1171.2858
$ws.Range("M110").Value2 = 421.25
$ws.Range("N110").Value2 = -5261.2858

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value2 = 1247.4722
$ws.Range("I136").Value2 = 837.43634
$ws.Range("J136").Value2 = 2574.0588
$ws.Range("K136").Value2 = 2512.30902
$ws.Range("L136").Value2 = 7722.176399999999
$ws.Range("M136").Value2 = 37.69098000000031
$ws.Range("N136").Value2 = -12822.1764

$ws = $wb.Worksheets.Item("BSM")
# Row 40: Can You Spare a Dolabra / Steel Dolabra
$ws.Range("H40").Value2 = 0
$ws.Range("J40").Value2 = 0
$ws.Range("L40").Value2 = 0
$ws.Range("N40").ClearContents()

# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value2 = 678.375
$ws.Range("J94").Value2 = 1453.3334
$ws.Range("L94").Value2 = 1453.3334
$ws.Range("N94").Value2 = -2355.3334

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value2 = 770.0625
$ws.Range("I107").Value2 = 728.7143
$ws.Range("J107").Value2 = 802.2222
$ws.Range("K107").Value2 = 728.7143
$ws.Range("L107").Value2 = 802.2222
$ws.Range("M107").Value2 = 1191.2857
$ws.Range("N107").Value2 = -4642.2222

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value2 = 412.2381
$ws.Range("I22").Value2 = 364.7143
$ws.Range("J22").Value2 = 507.2857
$ws.Range("K22").Value2 = 364.7143
$ws.Range("L22").Value2 = 507.2857
$ws.Range("M22").Value2 = -14.71429999999998
$ws.Range("N22").Value2 = -1207.2857

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value2 = 3088460
$ws.Range("I31").Value2 = 1413.0454
$ws.Range("J31").Value2 = 16671466
$ws.Range("K31").Value2 = 1413.0454
$ws.Range("L31").Value2 = 16671466
$ws.Range("M31").Value2 = -1118.0454
$ws.Range("N31").Value2 = -16672056

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value2 = 3088460
$ws.Range("I34").Value2 = 1413.0454
$ws.Range("J34").Value2 = 16671466
$ws.Range("K34").Value2 = 1413.0454
$ws.Range("L34").Value2 = 16671466
$ws.Range("M34").Value2 = -1211.0454
$ws.Range("N34").Value2 = -16671870

# Row 36: Toys of Summer / Steel Spear
$ws.Range("H36").Value2 = 10986.223
$ws.Range("I36").Value2 = 2032
$ws.Range("J36").Value2 = 15463.333
$ws.Range("K36").Value2 = 2032
$ws.Range("L36").Value2 = 15463.333
$ws.Range("M36").Value2 = -1644
$ws.Range("N36").Value2 = -16239.333

# Row 40: Ceremonial Spears / Steel Spear
$ws.Range("H40").Value2 = 10986.223
$ws.Range("I40").Value2 = 2032
$ws.Range("J40").Value2 = 15463.333
$ws.Range("K40").Value2 = 2032
$ws.Range("L40").Value2 = 15463.333
$ws.Range("M40").Value2 = -1872
$ws.Range("N40").Value2 = -15783.333

# Row 86: Birch, Please / Birch Lumber
$ws.Range("H86").Value2 = 90911020
$ws.Range("I86").Value2 = 142858900
$ws.Range("J86").Value2 = 2248.5
$ws.Range("K86").Value2 = 142858900
$ws.Range("L86").Value2 = 2248.5
$ws.Range("M86").Value2 = -142857777
$ws.Range("N86").Value2 = -4494.5

# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value2 = 90911020
$ws.Range("I89").Value2 = 142858900
$ws.Range("J89").Value2 = 2248.5
$ws.Range("K89").Value2 = 714294500
$ws.Range("L89").Value2 = 11242.5
$ws.Range("M89").Value2 = -714288884
$ws.Range("N89").Value2 = -22474.5

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value2 = 694.5
$ws.Range("I107").Value2 = 767.6111
$ws.Range("J107").Value2 = 621.3889
$ws.Range("K107").Value2 = 767.6111
$ws.Range("L107").Value2 = 621.3889
$ws.Range("M107").Value2 = 1152.3889
$ws.Range("N107").Value2 = -4461.3889

$ws = $wb.Worksheets.Item("CUL")
# Row 86: Let's Not Get Sappy / Birch Syrup
$ws.Range("H86").Value2 = 544.7692
$ws.Range("I86").Value2 = 556.5
$ws.Range("J86").Value2 = 534.7143
$ws.Range("K86").Value2 = 1669.5
$ws.Range("L86").Value2 = 1604.1429
$ws.Range("M86").Value2 = -483.5
$ws.Range("N86").Value2 = -3976.1429

# Row 89: Luxury Spillover (L) / Birch Syrup
$ws.Range("H89").Value2 = 544.7692
$ws.Range("I89").Value2 = 556.5
$ws.Range("J89").Value2 = 534.7143
$ws.Range("K89").Value2 = 5008.5
$ws.Range("L89").Value2 = 4812.428699999999
$ws.Range("M89").Value2 = 919.5
$ws.Range("N89").Value2 = -16668.4287

# Row 129: Comfort Food / Yakow Moussaka
$ws.Range("H129").Value2 = 1533.5385
$ws.Range("I129").Value2 = 1245.1666
$ws.Range("J129").Value2 = 1780.7142
$ws.Range("K129").Value2 = 3735.4998
$ws.Range("L129").Value2 = 5342.142599999999
$ws.Range("M129").Value2 = 1264.5002
$ws.Range("N129").Value2 = -15342.1426

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value2 = 6126094
$ws.Range("I131").Value2 = 22611.111
$ws.Range("J131").Value2 = 10049761
$ws.Range("K131").Value2 = 67833.333
$ws.Range("L131").Value2 = 30149283
$ws.Range("M131").Value2 = -62793.333
$ws.Range("N131").Value2 = -30159363

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value2 = 3128.9167
$ws.Range("I80").Value2 = 998.2
$ws.Range("J80").Value2 = 4650.857
$ws.Range("K80").Value2 = 998.2
$ws.Range("L80").Value2 = 4650.857
$ws.Range("M80").Value2 = -0.2000000000000455
$ws.Range("N80").Value2 = -6646.857

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value2 = 3128.9167
$ws.Range("I83").Value2 = 998.2
$ws.Range("J83").Value2 = 4650.857
$ws.Range("K83").Value2 = 4991
$ws.Range("L83").Value2 = 23254.285
$ws.Range("M83").Value2 = 1
$ws.Range("N83").Value2 = -33238.285

# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value2 = 427.4516
$ws.Range("I107").Value2 = 531.3158
$ws.Range("K107").Value2 = 531.3158
$ws.Range("M107").Value2 = 1388.6842

$ws = $wb.Worksheets.Item("LTW")
# Row 9: From the Sands to the Stage / Leather Himantes
$ws.Range("H9").Value2 = 677.375
$ws.Range("I9").Value2 = 383.8
$ws.Range("K9").Value2 = 383.8
$ws.Range("M9").Value2 = -159.8

# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value2 = 1246.4615
$ws.Range("I22").Value2 = 1440
$ws.Range("K22").Value2 = 1440
$ws.Range("M22").Value2 = -1145

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value2 = 1246.4615
$ws.Range("I27").Value2 = 1440
$ws.Range("K27").Value2 = 1440
$ws.Range("M27").Value2 = -1333

# Row 38: Emergency Patches / Skull Eyepatch
$ws.Range("H38").Value2 = 5000
$ws.Range("J38").Value2 = 5000
$ws.Range("L38").Value2 = 5000
$ws.Range("N38").Value2 = -5820

$ws = $wb.Worksheets.Item("WVR")
# Row 39: By the Short Hairs / Velveteen Robe
$ws.Range("H39").Value2 = 13022
$ws.Range("I39").Value2 = 13022
$ws.Range("K39").Value2 = 13022
$ws.Range("M39").Value2 = -12609

# Row 46: Crunching the Numbers / Linen Hat
$ws.Range("H46").Value2 = 75247.5
$ws.Range("J46").Value2 = 75247.5
$ws.Range("L46").Value2 = 75247.5
$ws.Range("N46").Value2 = -75709.5

# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value2 = 383.15384
$ws.Range("I107").Value2 = 320.1111
$ws.Range("J107").Value2 = 525
$ws.Range("K107").Value2 = 960.3333
$ws.Range("L107").Value2 = 1575
$ws.Range("M107").Value2 = 959.6667
$ws.Range("N107").Value2 = -5415

# Row 134: Cloth for Canvas / Mountain Linen
$ws.Range("H134").Value2 = 75247.5
$ws.Range("J134").Value2 = 75247.5
$ws.Range("L134").Value2 = 225742.5
$ws.Range("N134").Value2 = -230812.5

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value2 = 2165.4075
$ws.Range("I136").Value2 = 2378.1538
$ws.Range("K136").Value2 = 7134.4614
$ws.Range("M136").Value2 = -4584.4614
